# Weekly fruit/vegetable price update: two new daily price records were
# added to the "Berenjena" (eggplant) price log for Vega Central Mapocho de
# Santiago. They belong at the top of the data block (most recent date,
# 2022-08-25), so insert two new rows right after the header/existing data
# boundary (before the old row 241), pushing the rest of the table down by
# two rows (old row 241 -> new row 243, ..., old row 298 -> new row 300).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 241; this shifts all the
# existing data (previously rows 241:298) down to rows 243:300 and extends
# the used range to A1:R300.
$ws.Rows("241:242").Insert()

# --- New row 241 ---
$ws.Cells.Item(241, 1).Value = 9
$ws.Cells.Item(241, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(241, 3).Value = "Metropolitana"
$ws.Cells.Item(241, 4).Value = 44798
$ws.Cells.Item(241, 5).Value = 13
$ws.Cells.Item(241, 6).Value = 100112001
$ws.Cells.Item(241, 7).Value = "Berenjena"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 110
$ws.Cells.Item(241, 11).Value = 11000
$ws.Cells.Item(241, 12).Value = 12000
$ws.Cells.Item(241, 13).Value = 11545
$ws.Cells.Item(241, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(241, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(241, 16).Value = 289
$ws.Cells.Item(241, 17).Value = 40
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# --- New row 242 ---
$ws.Cells.Item(242, 1).Value = 9
$ws.Cells.Item(242, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(242, 3).Value = "Metropolitana"
$ws.Cells.Item(242, 4).Value = 44798
$ws.Cells.Item(242, 5).Value = 13
$ws.Cells.Item(242, 6).Value = 100112001
$ws.Cells.Item(242, 7).Value = "Berenjena"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 260
$ws.Cells.Item(242, 11).Value = 11000
$ws.Cells.Item(242, 12).Value = 12000
$ws.Cells.Item(242, 13).Value = 11865
$ws.Cells.Item(242, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(242, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(242, 16).Value = 237
$ws.Cells.Item(242, 17).Value = 50
$ws.Cells.Item(242, 18).Value = "Hortaliza"
